# "début ampl phase 1" - set up first data for the AMPL model on the
# "Fromage" sheet: a new "Gain" column next to the existing table, and a
# small matrix (rows 14-18) relating cheese types to lines/limits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fromage")

# New "Gain" header + values for the existing lot table (rows 1-7)
$ws.Range("F1").Value = "Gain"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 7
$ws.Range("F7").Value = 3

# New small matrix further down the sheet (rows 14-18)
$ws.Range("C14").Value = "Brie"
$ws.Range("D14").Value = "Sauvagine"
$ws.Range("E14").Value = "Dorémi"
$ws.Range("F14").Value = "Bleu"

$ws.Range("B15").Value = "Brie"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0

$ws.Range("B16").Value = "Sauv"
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0

$ws.Range("B17").Value = "Do"
$ws.Range("E17").Value = 0

$ws.Range("B18").Value = "Bleu"
$ws.Range("C18").Value = 30
$ws.Range("E18").Value = 0

# Match the author's final selection on the sheet
[void]$ws.Range("C2").Select()
